# Insert a new weekly price record at row 120 (pushing existing rows 120-215
# down to 121-216) for "Vega Monumental Concepción - Brócoli".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 120, shifting rows 120:215 down to 121:216
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new record
$ws.Cells.Item(120, 1).Value2 = 11
$ws.Cells.Item(120, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(120, 3).Value2 = "Bíobío"
$ws.Cells.Item(120, 4).Value2 = 44589
$ws.Cells.Item(120, 5).Value2 = 8
$ws.Cells.Item(120, 6).Value2 = 100112023
$ws.Cells.Item(120, 7).Value2 = "Brócoli"
$ws.Cells.Item(120, 8).Value2 = "Sin especificar"
$ws.Cells.Item(120, 9).Value2 = "Primera"
$ws.Cells.Item(120, 10).Value2 = 1300
$ws.Cells.Item(120, 11).Value2 = 700
$ws.Cells.Item(120, 12).Value2 = 750
$ws.Cells.Item(120, 13).Value2 = 723
$ws.Cells.Item(120, 14).Value2 = "`$/unidad"
$ws.Cells.Item(120, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(120, 16).Value2 = 723
$ws.Cells.Item(120, 17).Value2 = 1
$ws.Cells.Item(120, 18).Value2 = "Hortaliza"
